$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old data that lived in column F (rows 4-12)
$ws.Range("F4:F12").ClearContents()

# Header row
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Location"
$ws.Range("C1").Value = "Address"
$ws.Range("D1").Value = "Contact Number"

# Data rows
$data = @(
    @("Abid Computers",           "Saddar",            "Shop - 47, Technocity"),
    @("Alam Gamerz",               "Saddar",            "Shop - 48, Technocity"),
    @("Qazi and Kazi Brothers",    "University Avenue", "Baithak, Habib University"),
    @("HU Technologies",           "University Avenue", "Presidential Suite, Habib University"),
    @("DotShaheen Tech",           "Sir Syed",          "Shop C57, Sir Syed Road"),
    @("PowerfulDotOperator",       "Sir Syed",          "Shop C52, Sir Syed Road"),
    @("Altoo Tech Point",          "Saddar",            "Shop - 49, Technocity"),
    @("Bhatti LiquiTech",          "Saddar",            "Shop - 68, Technocity"),
    @("High Khalid Limited",       "Saddar",            "Shop - 69, Technocity")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}

# Column widths
$ws.Columns.Item(1).ColumnWidth = 24
$ws.Columns.Item(2).ColumnWidth = 20.140625

# Selection matches the authored state
$ws.Range("B13").Select()
